# [Kadastro App] Yeni kayit eklendi: 3013
# Appends the new "3013" record as row 73 on both the "Kayitlar" and the
# "Erdemli" sheets (the workbook keeps a master log plus a per-unit copy),
# matching how every other row in these sheets stores its values as text.

$wb = $excel.ActiveWorkbook

$newRowValues = @("3013", "2025-09-11", "Erdemli", "1", "3B", "SERDAR ARSLAN (Tekniker), ÖZKAN AKBAŞ (Mühendis)")
$targetRow = 73

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Leading apostrophes force text storage, same as the existing rows,
    # so numeric-looking values ("3013", "2025-09-11", "1") don't get
    # reinterpreted as a number/date.
    $ws.Cells.Item($targetRow, 1).Value = "'" + $newRowValues[0]
    $ws.Cells.Item($targetRow, 2).Value = "'" + $newRowValues[1]
    $ws.Cells.Item($targetRow, 3).Value = $newRowValues[2]
    $ws.Cells.Item($targetRow, 4).Value = "'" + $newRowValues[3]
    $ws.Cells.Item($targetRow, 5).Value = $newRowValues[4]
    $ws.Cells.Item($targetRow, 6).Value = $newRowValues[5]
}
